# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" bullets (under the "Impact"
# sub-heading) from job-duty style statements into short, impact-focused
# accomplishment statements, and drop two of the six bullets entirely.

$d = $word.ActiveDocument

# Locate the "Impact" sub-heading that sits under "KEY ACHIEVEMENTS AND IMPACT"
# (there are other "Impact:" paragraphs in KEY PROJECTS, so anchor on the
# Heading3-styled paragraph whose text is exactly "Impact").
$impactHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.Trim() -eq "Impact" -and $para.Style.NameLocal -eq "Heading 3") {
        $impactHeading = $para
        break
    }
}

$bullet1 = $impactHeading.Next()
$bullet2 = $bullet1.Next()
$bullet3 = $bullet2.Next()
$bullet4 = $bullet3.Next()
$bullet5 = $bullet4.Next()
$bullet6 = $bullet5.Next()

# 1) "Built redistricting platform ... 89 organizations" ->
#    "Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
$null = $bullet1.Range.Find.Execute(
    "Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%", 2)

# 2) "Designed ETL pipelines ..." -> "$4.7M savings enabled nonprofit access"
$null = $bullet2.Range.Find.Execute(
    "Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "`$4.7M savings enabled nonprofit access", 2)

# 3) "Trigonometric algorithm for boundary estimation ..." ->
#    "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
$null = $bullet3.Range.Find.Execute(
    "Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations", 2)

# 6) "Built cloud-based data warehouse solutions ..." -> "Real-time collaboration at national scale"
$null = $bullet6.Range.Find.Execute(
    "Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Real-time collaboration at national scale", 2)

# 4) and 5) are dropped entirely: "Discovered systematic race coding errors ..." and
#    "Achieved 87% prediction accuracy ..." - delete highest index first so the
#    other paragraph reference stays valid.
$bullet5.Range.Delete()
$bullet4.Range.Delete()
